$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value (all columns D/E store formatted text strings,
# so we force a Text number format before assignment to avoid Excel re-parsing the
# string as a number/percentage and losing the exact displayed precision, then reset
# the style back to Normal so no stray formatting is left behind on the cell.
$updates = @{
    "D2" = "261.32"
    "E2" = "1.65%"
    "D3" = "27.37"
    "E3" = "1.29%"
    "D4" = "4.719"
    "E4" = "0.70%"
    "E5" = "3.13%"
    "D6" = "6.677"
    "E6" = "1.03%"
    "D7" = "0.8463"
    "E7" = "-0.44%"
    "D8" = "0.9219"
    "E8" = "-0.65%"
    "D9" = "0.1404"
    "E9" = "1.96%"
    "D10" = "0.05011"
    "E10" = "17.32%"
    "D11" = "0.07133"
    "E11" = "1.37%"
    "D12" = "0.03115"
    "E12" = "1.97%"
    "D13" = "0.09069"
    "E13" = "-0.41%"
    "D14" = "0.001541"
    "E14" = "-0.06%"
    "D15" = "0.0006076"
    "E15" = "0.13%"
    "D16" = "0.006142"
    "E16" = "2.14%"
    "D17" = "3.451"
    "E17" = "-0.52%"
    "E18" = "-0.88%"
    "D19" = "2.186"
    "E19" = "-1.15%"
    "D20" = "0.3126"
    "E20" = "1.47%"
    "D21" = "0.1305"
    "E21" = "0.84%"
    "D22" = "4.093"
    "E22" = "4.94%"
    "D23" = "0.04245"
    "E23" = "-0.06%"
    "E24" = "0.27%"
    "E25" = "-9.00%"
    "D26" = "0.0001201"
    "E26" = "0.10%"
    "D27" = "0.0001576"
    "E27" = "3.46%"
    "E40" = "1.95%"
    "D41" = "0.1115"
    "E41" = "1.45%"
    "D42" = "0.004090"
    "E42" = "-34.22%"
    "D43" = "0.01637"
    "E43" = "21.79%"
    "D44" = "0.002218"
    "E44" = "0.83%"
    "D45" = "0.00005262"
    "E45" = "-1.59%"
    "E47" = "26.38%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
